$d = $word.ActiveDocument

function Set-ParaXml($para, [string]$innerXml) {
    $rngTarget = $para.Range
    $bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rngTarget.InsertXML($pkgXml, "Replace")
}

# Locate the bookmark ("_GoBack") paragraph, which is the final paragraph of the document,
# and the empty paragraph that immediately precedes it.
$bm = $d.Bookmarks.Item("_GoBack")
$bookmarkParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $bm.Range.Start) {
        $bookmarkParaIndex = $i
    }
}
$precedingParaIndex = $bookmarkParaIndex - 1

# Remove the now-superfluous empty paragraph that precedes the bookmark paragraph;
# its blank content will be recreated (without bold) after the bookmark paragraph below.
$precedingPara = $d.Paragraphs.Item($precedingParaIndex)
$precedingPara.Range.Delete()

# Re-fetch the bookmark paragraph (its index shifted down by one after the deletion above).
$bookmarkParaIndex = $bookmarkParaIndex - 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# Insert the new content paragraphs immediately before the bookmark paragraph.
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()

# Fill in the text/formatting of each newly created paragraph, in document order.
$firstNewIndex = $bookmarkParaIndex
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 0)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>DATEADD</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 1)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Se utiliza para sumar una cantidad de tiempo a una fecha</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 2)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>DATEDIFF</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 3)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Sirve para calcular la diferencia entre dos fechas</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 4)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">FUNCIONES DE CADENA </w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 5)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>CONCAT</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 6)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Sirve para unir cadenas de texto</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 7)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>LEN</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 8)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Nos muestra la cantidad de caracteres que existen</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 9)) '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>SUBSTRING</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 10)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Extrae una parte del texto o carácter </w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 11)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>UPPER.-</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Convierte el texto a mayúscula </w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 12)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">LOWER. – </w:t></w:r><w:r><w:t xml:space="preserve">Parecido al </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>upper</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> pero convierte en minúsculas</w:t></w:r>'
Set-ParaXml ($d.Paragraphs.Item($firstNewIndex + 13)) '<w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">LTRIM o RTRIM. – </w:t></w:r><w:r><w:t xml:space="preserve">Quita espacios en blanco al inicio y al final </w:t></w:r>'

# Re-fetch the bookmark paragraph (index shifted down by the inserted paragraphs) and
# strip the bold paragraph-mark formatting, keeping its bookmark intact.
$bookmarkParaIndex = $firstNewIndex + 14
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
Set-ParaXml $bookmarkPara '<w:pPr><w:jc w:val="both"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

# Append the two trailing empty paragraphs (no bold, then bold) after the bookmark paragraph.
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.InsertParagraphAfter()
$emptyPara1 = $d.Paragraphs.Item($bookmarkParaIndex + 1)
Set-ParaXml $emptyPara1 '<w:pPr><w:jc w:val="both"/></w:pPr>'
$emptyPara1 = $d.Paragraphs.Item($bookmarkParaIndex + 1)
$emptyPara1.Range.InsertParagraphAfter()
$emptyPara2 = $d.Paragraphs.Item($bookmarkParaIndex + 2)
Set-ParaXml $emptyPara2 '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr>'

Write-Host "Final paragraph count:" $d.Paragraphs.Count
